# Correction des données générées pour les fichiers JSON :
# l'année des dates (colonne A) passe de 2015 à 2017, et les moyennes
# (colonne E) sont recalculées pour les étudiants concernés (certaines
# lignes conservent leur moyenne d'origine lorsque la donnée n'existait
# pas / n'a pas changé).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PPP")

# Table: numéro de ligne, nouvelle date (colonne A), nouvelle moyenne (colonne E, $null = inchangé)
$updates = @(
    @(3, 20170926, 19),
    @(4, 20170927, 5),
    @(5, 20170928, 12),
    @(6, 20170929, 16),
    @(7, 20170930, 20),
    @(8, 20170931, 12),
    @(9, 20170932, 9),
    @(10, 20170933, 9),
    @(11, 20170934, 5),
    @(12, 20170935, 14),
    @(13, 20170936, 9),
    @(14, 20170937, 15),
    @(15, 20170938, 19),
    @(16, 20170939, 20),
    @(17, 20170940, 16),
    @(18, 20170941, 13),
    @(19, 20170942, 5),
    @(20, 20170943, $null),
    @(21, 20170944, 5),
    @(22, 20170945, 12),
    @(23, 20170946, $null),
    @(24, 20170947, 17),
    @(25, 20170948, 5),
    @(26, 20170949, 9),
    @(27, 20170950, 9),
    @(28, 20170951, 7),
    @(29, 20170952, 9),
    @(30, 20170953, 15),
    @(31, 20170954, 6),
    @(32, 20170955, $null),
    @(33, 20170956, 8),
    @(34, 20170957, 15),
    @(35, 20170958, 14),
    @(36, 20170959, 8),
    @(37, 20170960, 12),
    @(38, 20170961, 20),
    @(39, 20170962, 10),
    @(40, 20170963, 12),
    @(41, 20170964, 11),
    @(42, 20170965, 6),
    @(43, 20170966, 14),
    @(44, 20170967, $null),
    @(45, 20170968, 6),
    @(46, 20170969, 14),
    @(47, 20170970, $null),
    @(48, 20170971, 12),
    @(49, 20170972, 17),
    @(50, 20170973, 7),
    @(51, 20170974, 15),
    @(52, 20170975, 16),
    @(53, 20170976, 20),
    @(54, 20170977, 14),
    @(55, 20170978, 10),
    @(56, 20170979, 16),
    @(57, 20170980, 11),
    @(58, 20170981, 10),
    @(59, 20170982, 7),
    @(60, 20170983, 15),
    @(61, 20170984, 12),
    @(62, 20170985, 18),
    @(63, 20170986, 8)
)

foreach ($update in $updates) {
    $row = $update[0]
    $newDate = $update[1]
    $newAvg = $update[2]

    $ws.Cells.Item($row, 1).Value = $newDate
    if ($null -ne $newAvg) {
        $ws.Cells.Item($row, 5).Value = $newAvg
    }
}
